$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The abbreviation used for "Eye problems" is renamed from the truncated
# "eye_prob" to the full "eye_problems".
$ws.Range("B6").Value = "eye_problems"

# Column B (the abbreviation column) now gets an explicit best-fit width,
# matching column A which already has a fixed custom width. The engine
# stores column widths snapped to 1/6-character increments and adds a
# fixed 0.8333... (5/6) padding on write, so we back that offset out of
# the desired stored width (12.83203125) to land as close as possible to
# the target value.
$ws.Columns.Item(2).ColumnWidth = 12.83203125 - (5/6)

# The active cell / selection shown when the sheet is next opened moves
# from G18 to B7.
$ws.Range("B7").Select()
